# "Generate Report for Handback"
# The handback for the de-de target succeeded (its version check passed), so the
# workbook's status/report cells move from "Ready for handoff" to
# "Handed back: in sync with en-US", the handoff/handback timestamps advance, and
# the (now resolved) version-mismatch error detail is cleared. Column widths for
# the Status / Error Detail columns are widened/narrowed to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: per-locale status columns (E = zh-cn, F = de-de) ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn detail sheet ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-09-01 18:58:36"
$wsZhCn.Range("P2").Value = ""

# --- de-de detail sheet ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-09-01 18:58:45"
$wsDeDe.Range("P2").Value = ""

# --- Column width adjustments (Status / Error Detail columns widen & narrow
#     to fit the new, differently-sized, text) ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
